$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2926
$ws.Range("F4").Value = 2926
$ws.Range("F5").Value = 6405
$ws.Range("F6").Value = 2511
$ws.Range("F8").Value = 414
$ws.Range("F9").Value = 51
$ws.Range("F11").Value = 2922
$ws.Range("F14").Value = 7309
$ws.Range("F15").Value = 329
$ws.Range("F16").Value = 63
$ws.Range("F17").Value = 105
$ws.Range("F18").Value = 236
$ws.Range("F20").Value = 481
$ws.Range("F21").Value = 8857
$ws.Range("F23").Value = 41
$ws.Range("F24").Value = 266
$ws.Range("F27").Value = 25
$ws.Range("F35").Value = 2613
$ws.Range("F38").Value = 38
$ws.Range("F40").Value = 728
$ws.Range("F41").Value = 3836
$ws.Range("F43").Value = 195
$ws.Range("F44").Value = 29
$ws.Range("F46").Value = 9
$ws.Range("F47").Value = 223
$ws.Range("F48").Value = 5
$ws.Range("F49").Value = 48
$ws.Range("F50").Value = 38

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 258

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2926
$ws.Range("F7").Value = 258
$ws.Range("F9").Value = 6405
$ws.Range("F12").Value = 51
$ws.Range("F19").Value = 7309
$ws.Range("F20").Value = 329
$ws.Range("F22").Value = 105
$ws.Range("F23").Value = 236
$ws.Range("F25").Value = 8857
$ws.Range("F27").Value = 41
$ws.Range("F30").Value = 25
$ws.Range("F38").Value = 2613
$ws.Range("F41").Value = 38
$ws.Range("F42").Value = 728
$ws.Range("F43").Value = 3836
$ws.Range("F44").Value = 195
$ws.Range("F48").Value = 223
$ws.Range("F49").Value = 48
$ws.Range("F50").Value = 38
